# Weekly update: a new daily price record was inserted ahead of the
# existing row 187, pushing every subsequent record down by one row
# (187-302 -> 188-303) and growing the table by one row (now A1:T303).
#
# Insert a fresh row at 187 (shifts 187:302 down to 188:303, carrying all
# formatting/styles with it) and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("187:187").Insert()

$ws.Cells.Item(187, 1).Value  = 4
$ws.Cells.Item(187, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(187, 3).Value  = "Los Lagos"
$ws.Cells.Item(187, 4).Value  = 44762
$ws.Cells.Item(187, 5).Value  = 10
$ws.Cells.Item(187, 6).Value  = "Fruta"
$ws.Cells.Item(187, 7).Value  = 100101
$ws.Cells.Item(187, 8).Value  = "Berries"
$ws.Cells.Item(187, 9).Value  = 100101007
$ws.Cells.Item(187, 10).Value = "Kiwi"
$ws.Cells.Item(187, 11).Value = "Hayward"
$ws.Cells.Item(187, 12).Value = "Primera"
$ws.Cells.Item(187, 13).Value = 200
$ws.Cells.Item(187, 14).Value = 12000
$ws.Cells.Item(187, 15).Value = 14000
$ws.Cells.Item(187, 16).Value = 13000
$ws.Cells.Item(187, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(187, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(187, 19).Value = 867
$ws.Cells.Item(187, 20).Value = 15
